$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style the "Other services:" / "Tasks:" / "Lista funkcjonalnosci:" headers
# (L2, N2, P2) so they match the same header look used elsewhere (e.g. B15),
# collapsing the now-redundant style onto the shared one.
$ws.Range("B15").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New "Wykorzystane technologie" (Technologies used) section, column P
$ws.Range("B15").Copy() | Out-Null
$ws.Range("P14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("P14").Value = "Wykorzystane technologie"

$ws.Range("P15").Value = "Entity Framework"
$ws.Range("P16").Value = "MSSQL"
$ws.Range("P17").Value = ".NET WPF"
$ws.Range("P18").Value = "SendGrid (EmailSender)"
$ws.Range("P19").Value = "PdfSharp (PdfGenerator)"
$ws.Range("P20").Value = "AutoMapper"

# --- New "Wzorce projektowe i architektura" (Design patterns and architecture) section
# (values entered in their original authoring order so the shared-string
# table comes out in the same sequence)
$ws.Range("P23").Value = "Dependency Injection"

$ws.Range("B15").Copy() | Out-Null
$ws.Range("P22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("P22").Value = "Wzorce projektowe i architektura"

$ws.Range("P25").Value = "Repository"
$ws.Range("P24").Value = "Factory"
$ws.Range("P26").Value = "Microservices"
$ws.Range("P27").Value = "MVVM"

# --- Move the active selection (cosmetic, matches the saved cursor position)
$ws.Range("F4").Select() | Out-Null
